$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D values are plain-text price strings that often look like numbers
# (e.g. "237.02", "1.000", "15.00"). Excel auto-converts such assignments to numeric
# cells (losing exact formatting / trailing zeros). To preserve them as literal text
# (matching the original inlineStr cells), we force the cell format to Text ("@")
# before assigning, then restore the default "Normal" style so no stray number format
# is left behind on the cell.

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.768.70'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.31%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.752.66'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -4.29%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '237.02'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.93%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.06%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5067'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -3.54%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '41.41'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -6.64%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2655'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.47%  '

# Row 10
$ws.Range('E10').Value = '  -9.81%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.752.54'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.44%  '

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.06904'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.57%  '

# Row 13
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '15.61'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -5.51%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6023'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -12.42%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.495'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -7.23%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '77.20'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -10.23%  '

# Row 17
$ws.Range('E17').Value = '  -0.09%  '

# Row 18
$ws.Range('E18').Value = '  -0.09%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '25.781.93'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.40%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000006842'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -6.47%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.69'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -11.26%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.979.93'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.03%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.084'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -9.06%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.239'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -7.93%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.201'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -10.31%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '137.40'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.51%  '

# Row 27
$ws.Range('E27').Value = '  -12.08%  '

# Row 28
$ws.Range('E28').Value = '  -9.48%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -9.15%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '102.64'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.82%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08199'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -5.89%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.684'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -9.17%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.449'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -10.33%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04512'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.52%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.000'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.03%  '

# Row 36
$ws.Range('E36').Value = '  -7.68%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9952'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -9.80%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.6057'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -13.61%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.698'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -11.67%  '

# Row 40
$ws.Range('E40').Value = '  -4.87%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.943'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -10.49%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.001'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.01%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '103.46'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.40%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3813'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -14.24%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.7401'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -13.85%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.916'
$ws.Range('D46').Style = "Normal"

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05466'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.75%  '

# Row 48
$ws.Range('E48').Value = '  -5.93%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.973'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -14.33%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.683'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -10.92%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '29.93'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -10.13%  '
